$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns AC and AD
$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

# New data values for rows 2 and 3
$ws.Range("AC2").Value = 0.8038585209003215
$ws.Range("AD2").Value = 0.89320388349514568

$ws.Range("AC3").Value = 0.80434782608695654
$ws.Range("AD3").Value = 0.89078498293515362

# Expand the selection to match the new used range
$ws.Range("A1:AD3").Select()
